$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 15 (pushes the SOFT BODY..COMPOUND SHAPES block down by one,
# old row 15 becomes row 16, ..., old row 22 becomes row 23). Insert() copies the
# formatting of the row above (row 14), which already matches the target styling
# for column F (s=5).
$ws.Rows.Item(15).Insert()

# New "Networks" table header (K1:M1), mirroring the Physics table header (F1:H1).
$ws.Range("K1").Value = "Networks"
$ws.Range("L1").Value = "Notes"
$ws.Range("M1").Value = "Status"
$ws.Range("F1:H1").Copy()
$ws.Range("K1:M1").PasteSpecial(-4122)

# New task row (row 15): a physics-derived TODO item moved into the Networks list.
$ws.Range("F15").Value = "Create an outofbound check"
$ws.Range("G15").Value = "Rermoves need for expensive floor/ wall collision check"
$ws.Range("H15").Value = "TODO"

# H15 needs the plain "TODO" (red) style used elsewhere (e.g. H2), not the style
# inherited from the row-insert (which copied H14's "UNDERWAY" style).
$ws.Range("H2").Copy()
$ws.Range("H15").PasteSpecial(-4122)

# New column K width (as wide as the Notes column needs to show the new text).
$ws.Columns.Item(11).ColumnWidth = 24

# Matches the selection captured in the saved workbook.
[void]$ws.Range("K19").Select()

[void]$excel.CutCopyMode
